$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Domantas Sabonis -> Davion Mitchell
$ws.Range("B3").Value = 15
$ws.Range("C3").Value = "Davion Mitchell"
$ws.Range("D3").Value = "PG"
$ws.Range("E3").Value = "6-2"
$ws.Range("F3").Value = 205
$ws.Range("G3").Value = "September 5, 1998"
$ws.Range("I3").Value = "1"
$ws.Range("J3").Value = "Auburn, Baylor"
$ws.Range("K3").Value = "https://www.basketball-reference.com/players/m/mitchda01.html"

# Row 4: Malik Monk -> Domantas Sabonis
$ws.Range("B4").Value = 10
$ws.Range("C4").Value = "Domantas Sabonis"
$ws.Range("D4").Value = "C"
$ws.Range("E4").Value = "6-11"
$ws.Range("F4").Value = 240
$ws.Range("G4").Value = "May 3, 1996"
$ws.Range("I4").Value = "6"
$ws.Range("J4").Value = "Gonzaga"
$ws.Range("K4").Value = "https://www.basketball-reference.com/players/s/sabondo01.html"

# Row 5: Davion Mitchell -> Malik Monk
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = "Malik Monk"
$ws.Range("D5").Value = "SG"
$ws.Range("E5").Value = "6-3"
$ws.Range("F5").Value = 200
$ws.Range("G5").Value = "February 4, 1998"
$ws.Range("I5").Value = "5"
$ws.Range("J5").Value = "Kentucky"
$ws.Range("K5").Value = "https://www.basketball-reference.com/players/m/monkma01.html"

# Row 7: De'Aaron Fox -> Kevin Huerter
$ws.Range("B7").Value = 9
$ws.Range("C7").Value = "Kevin Huerter"
$ws.Range("D7").Value = "SG"
$ws.Range("E7").Value = "6-7"
$ws.Range("F7").Value = 190
$ws.Range("G7").Value = "August 27, 1998"
$ws.Range("I7").Value = "4"
$ws.Range("J7").Value = "Maryland"
$ws.Range("K7").Value = "https://www.basketball-reference.com/players/h/huertke01.html"

# Row 8: Kevin Huerter -> De'Aaron Fox
$ws.Range("B8").Value = 5
$ws.Range("C8").Value = "De'Aaron Fox"
$ws.Range("D8").Value = "PG"
$ws.Range("E8").Value = "6-3"
$ws.Range("F8").Value = 185
$ws.Range("G8").Value = "December 20, 1997"
$ws.Range("I8").Value = "5"
$ws.Range("J8").Value = "Kentucky"
$ws.Range("K8").Value = "https://www.basketball-reference.com/players/f/foxde01.html"

# Row 18: PJ Dozier -> Kessler Edwards (no jersey number)
$ws.Range("B18").ClearContents()
$ws.Range("C18").Value = "Kessler Edwards"
$ws.Range("D18").Value = "SF"
$ws.Range("E18").Value = "6-8"
$ws.Range("F18").Value = 215
$ws.Range("G18").Value = "August 9, 2000"
$ws.Range("I18").Value = "1"
$ws.Range("J18").Value = "Pepperdine"
$ws.Range("K18").Value = "https://www.basketball-reference.com/players/e/edwarke02.html"
